$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Row 2 updates - new iProctor smoke testcase data
$ws.Range("A2").Value = "jMtSL985"
$ws.Range("B2").Value = 23110301
$ws.Range("C2").Value = "ipwxppd82"
$ws.Range("D2").Value = "kR82P&%u"
$ws.Range("F2").Value = "OMCJfClL"
$ws.Range("G2").Value = "CELt"
